$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.432.33"
$ws.Range("E2").Value = "  +0.23%  "
$ws.Range("D3").Value = "'1.866.50"
$ws.Range("E3").Value = "  -0.70%  "
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("B5").Value = "XRP"
$ws.Range("C5").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$ws.Range("D5").Value = "'0.7072"
$ws.Range("E5").Value = "  -0.80%  "
$ws.Range("B6").Value = "BNB"
$ws.Range("C6").Value = "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"
$ws.Range("D6").Value = "'243.40"
$ws.Range("E6").Value = "  +0.28%  "
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("D8").Value = "'0.3141"
$ws.Range("E8").Value = "  -1.14%  "
$ws.Range("D9").Value = "'0.07856"
$ws.Range("E9").Value = "  -2.34%  "
$ws.Range("D10").Value = "'24.48"
$ws.Range("E10").Value = "  -2.65%  "
$ws.Range("D11").Value = "'0.08025"
$ws.Range("E11").Value = "  -3.85%  "
$ws.Range("D12").Value = "'1.873.49"
$ws.Range("E12").Value = "  -0.75%  "
$ws.Range("D13").Value = "'5.202"
$ws.Range("E13").Value = "  -1.24%  "
$ws.Range("D14").Value = "'93.42"
$ws.Range("E14").Value = "  -1.50%  "
$ws.Range("D15").Value = "'0.7010"
$ws.Range("E15").Value = "  -2.44%  "
$ws.Range("D16").Value = "'6.462"
$ws.Range("E16").Value = "  +1.37%  "
$ws.Range("D17").Value = "'29.449.34"
$ws.Range("E17").Value = "  +0.25%  "
$ws.Range("D18").Value = "'0.000008366"
$ws.Range("E18").Value = "  -3.25%  "
$ws.Range("D19").Value = "'252.83"
$ws.Range("E19").Value = "  +3.90%  "
$ws.Range("D20").Value = "'2.132.85"
$ws.Range("E20").Value = "  -0.79%  "
$ws.Range("D21").Value = "'13.14"
$ws.Range("E21").Value = "  -1.46%  "
$ws.Range("E22").Value = "  -0.08%  "
$ws.Range("D23").Value = "'7.611"
$ws.Range("E23").Value = "  -2.82%  "
$ws.Range("E24").Value = "  -0.08%  "
$ws.Range("D25").Value = "'0.1557"
$ws.Range("E25").Value = "  -1.22%  "
$ws.Range("D26").Value = "'9.019"
$ws.Range("E26").Value = "  -0.90%  "
$ws.Range("D27").Value = "'160.79"
$ws.Range("E27").Value = "  -1.51%  "
$ws.Range("E28").Value = "  +0.57%  "
$ws.Range("D29").Value = "'1.500"
$ws.Range("E29").Value = "  -0.43%  "
$ws.Range("D30").Value = "'4.324"
$ws.Range("E30").Value = "  -2.61%  "
$ws.Range("D31").Value = "'4.282"
$ws.Range("E31").Value = "  -1.69%  "
$ws.Range("D32").Value = "'1.210"
$ws.Range("E32").Value = "  +0.28%  "
$ws.Range("D33").Value = "'0.05306"
$ws.Range("E33").Value = "  -2.14%  "
$ws.Range("D34").Value = "'1.885"
$ws.Range("E34").Value = "  -3.35%  "
$ws.Range("D35").Value = "'0.7518"
$ws.Range("E35").Value = "  -2.78%  "
$ws.Range("D36").Value = "'1.168"
$ws.Range("E36").Value = "  -1.86%  "
$ws.Range("D37").Value = "'2.713"
$ws.Range("E37").Value = "  +1.06%  "
$ws.Range("D38").Value = "'0.01879"
$ws.Range("E38").Value = "  -0.72%  "
$ws.Range("D39").Value = "'1.262.28"
$ws.Range("E39").Value = "  -0.70%  "
$ws.Range("D40").Value = "'2.741"
$ws.Range("E40").Value = "  -0.41%  "
$ws.Range("D41").Value = "'0.8973"
$ws.Range("E41").Value = "  -1.47%  "
$ws.Range("D42").Value = "'108.89"
$ws.Range("E42").Value = "  -4.17%  "
$ws.Range("D43").Value = "'5.959"
$ws.Range("E43").Value = "  -8.56%  "
$ws.Range("D44").Value = "'71.32"
$ws.Range("E44").Value = "  -4.53%  "
$ws.Range("E45").Value = "  -0.04%  "
$ws.Range("E46").Value = "  -0.55%  "
$ws.Range("D47").Value = "'2.035.92"
$ws.Range("E47").Value = "  -0.14%  "
$ws.Range("B48").Value = "RenderToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D48").Value = "'1.789"
$ws.Range("E48").Value = "  -1.28%  "
$ws.Range("B49").Value = "Mantle"
$ws.Range("C49").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D49").Value = "'0.5177"
$ws.Range("E49").Value = "  -0.90%  "
$ws.Range("D50").Value = "'9.533"
$ws.Range("E50").Value = "  -0.36%  "
$ws.Range("D51").Value = "'0.4309"
$ws.Range("E51").Value = "  -1.83%  "
